# Adding MPA test automation upload file
# Update the generated batch/voucher reference numbers and their
# corresponding line-count values on the "Data" sheet.
#
# Mapping applied to rows 6-29:
#   column K: 60000231 -> 60000288
#   column L: 179      -> 229
#   column N: 60000232 -> 60000289
#   column O: 180      -> 230

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("K6").Value = 60000288
$ws.Range("L7").Value = 229
$ws.Range("N7").Value = 60000289
$ws.Range("L8").Value = 229
$ws.Range("O8").Value = 230
$ws.Range("K9").Value = 60000288
$ws.Range("L10").Value = 229
$ws.Range("K11").Value = 60000288
$ws.Range("N11").Value = 60000289
$ws.Range("L12").Value = 229
$ws.Range("N12").Value = 60000289
$ws.Range("L13").Value = 229
$ws.Range("O13").Value = 230
$ws.Range("K14").Value = 60000288
$ws.Range("L15").Value = 229
$ws.Range("K16").Value = 60000288
$ws.Range("N16").Value = 60000289
$ws.Range("L17").Value = 229
$ws.Range("N17").Value = 60000289
$ws.Range("L18").Value = 229
$ws.Range("O18").Value = 230
$ws.Range("K19").Value = 60000288
$ws.Range("L20").Value = 229
$ws.Range("K21").Value = 60000288
$ws.Range("N21").Value = 60000289
$ws.Range("L22").Value = 229
$ws.Range("N22").Value = 60000289
$ws.Range("L23").Value = 229
$ws.Range("O23").Value = 230
$ws.Range("K24").Value = 60000288
$ws.Range("L25").Value = 229
$ws.Range("K26").Value = 60000288
$ws.Range("N26").Value = 60000289
$ws.Range("L27").Value = 229
$ws.Range("N27").Value = 60000289
$ws.Range("L28").Value = 229
$ws.Range("O28").Value = 230
$ws.Range("K29").Value = 60000288
